# Update the "想去人数" (attendance interest count) figures on the
# "展览" and "全部类型" worksheets to reflect the newly generated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 303
    "F5"  = 30
    "F7"  = 10254
    "F11" = 129
    "F12" = 12
    "F14" = 32
    "F18" = 803
    "F20" = 97
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
